$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the exam note in B2 to append the extra air-temperature reading
$ws.Range("B2").Value = "``Date of Exam: 12/25/1999. Air temperature: 79.0 degF. Air pressure: 1013.2 mBar. Air temperature: 81.0 degF``"

# Move the active cell selection from B3 to B2
$ws.Range("B2").Select()
